$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new "Data Structures" worksheet after "Datacamp Courses" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Data Structures"

# --- Fill in the topic list (rows 2-6 were entered before the header row) ---
$ws2.Range("A2").Value = "Arrays in Data Structures"
$ws2.Range("A3").Value = "Array Operations"
$ws2.Range("A4").Value = "Array Operations Continued"
$ws2.Range("A5").Value = "Pointers & Arrays"
$ws2.Range("A6").Value = "2D Arrays"

# --- Header row, added after the first few topic rows ---
$ws2.Range("A1").Value = "Jenny's Videos"
$ws2.Range("B1").Value = "Completed"
$ws2.Range("A1:B1").Font.Bold = $true

# --- Remaining topic rows ---
$ws2.Range("A7").Value = "Pointers in 2D Arrays"
$ws2.Range("A8").Value = "Linked Lists"
$ws2.Range("A9").Value = "Types of Linked Lists"
$ws2.Range("A10").Value = "Arrays vs Linked Lists"
$ws2.Range("A11").Value = "Linked Lists Implementations"
$ws2.Range("A12").Value = "Insert a Node in a Singly LL"
$ws2.Range("A13").Value = "Delete a Node in a Singly LL"
$ws2.Range("A14").Value = "Length of LL"
$ws2.Range("A15").Value = "Reverse a Linked List"
$ws2.Range("A16").Value = "Doubly LL"
$ws2.Range("A17").Value = "Implement a Doubly LL"
$ws2.Range("A18").Value = "Insertion in Doubly LL"
$ws2.Range("A19").Value = "Reverse a Doubly LL"
$ws2.Range("A20").Value = "Circular LL"
$ws2.Range("A21").Value = "Implementation of Circular LL"
$ws2.Range("A22").Value = "CLL insertion"
$ws2.Range("A23").Value = "CLL deletion"
$ws2.Range("A24").Value = "Reverse a CLL"
$ws2.Range("A25").Value = "Stacks"
$ws2.Range("A26").Value = "Implementation of Stacks using Arrays"
$ws2.Range("A27").Value = "Implementation of Stacks using Linked Lists"
$ws2.Range("A28").Value = "Queues"
$ws2.Range("A29").Value = "Implementation of Queues using arrays"
$ws2.Range("A30").Value = "Implementation of Queues using LL"
$ws2.Range("A31").Value = "Circular Queue using Arrays"
$ws2.Range("A32").Value = "Circular Queue using LL"
$ws2.Range("A33").Value = "Implementation of Queues using Stack"
$ws2.Range("A34").Value = "Deque "
$ws2.Range("A35").Value = "Implementation of Deque using circular queue"
$ws2.Range("A36").Value = "Implementation of Deque using circular array"
$ws2.Range("A37").Value = "Trees"
$ws2.Range("A38").Value = "Binary Trees and its types"
$ws2.Range("A39").Value = "Binary Tree Implementation"
$ws2.Range("A40").Value = "Binary Tree using Arrays"
$ws2.Range("A41").Value = "Binary Tree traversal"
$ws2.Range("A42").Value = "Binary Search Trees"
$ws2.Range("A43").Value = "AVL Trees"
$ws2.Range("A44").Value = "Red Black Tree"
$ws2.Range("A45").Value = "Splay Trees"
$ws2.Range("A46").Value = "B-Trees"
$ws2.Range("A47").Value = "Graphs"
$ws2.Range("A48").Value = "Search Algorithms"
$ws2.Range("A49").Value = "Sorting Algorithms"
$ws2.Range("A50").Value = "Hashing Techniques"

# --- Column sizing / view state for the new sheet ---
$ws2.Columns.Item(1).ColumnWidth = 37.6640625
$ws2.Activate()
$excel.ActiveWindow.Zoom = 161
$ws2.Range("E39").Select()

# --- Record the just-completed Datacamp course on the first sheet ---
$ws1.Range("A29").Value = "Dealing with Missing Data in Python"
$ws1.Range("B29").Value = 44176
$ws1.Range("B29").NumberFormat = "mm-dd-yy"

# --- Switch back to the first sheet and leave the selection on the new row ---
$ws1.Activate()
$ws1.Range("C29").Select()
